# Method of Hierarchy Analysis
# - Narrow the comparison-matrix columns (car model names now wrap onto
#   multiple lines instead of running the sheet wide).
# - Break the long car-model names / header captions onto multiple lines
#   (space -> line break) so they read well in the now-narrower columns.
# - Refresh two of the computed priority-vector values (row 4 / row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -------------------------------------------------
# ColumnWidth is in "characters"; the saved OOXML <col width> ends up as
# ColumnWidth + 5/6. Back the character width out from the target OOXML
# width so the saved file lands as close as possible to it.
$ws.Columns.Item(1).ColumnWidth = 14.4 - 5/6
$ws.Columns.Item(2).ColumnWidth = 6 - 5/6
$ws.Columns.Item(3).ColumnWidth = 14.4 - 5/6
$ws.Columns.Item(4).ColumnWidth = 10.8 - 5/6
$ws.Columns.Item(5).ColumnWidth = 10.8 - 5/6
$ws.Columns.Item(6).ColumnWidth = 9.6 - 5/6
$ws.Columns.Item(7).ColumnWidth = 10.8 - 5/6
$ws.Columns.Item(8).ColumnWidth = 15.6 - 5/6

# --- Header row (row 1): car names split across two (or three) lines ---
$ws.Range("B1").Value = "Kia" + [char]10 + "Rio"
$ws.Range("C1").Value = "Volkswagen" + [char]10 + "Golf"
$ws.Range("D1").Value = "Toyota" + [char]10 + "Corolla"
$ws.Range("E1").Value = "Skoda" + [char]10 + "Octavia"
$ws.Range("F1").Value = "BMW" + [char]10 + "3" + [char]10 + "Series"
$ws.Range("G1").Value = "Hyundai" + [char]10 + "Solaris"
$ws.Range("H1").Value = "Вектор" + [char]10 + "приоритетов"

# --- Row labels (column A, rows 2-7): same car names, same line breaks ---
$ws.Range("A2").Value = "Kia" + [char]10 + "Rio"
$ws.Range("A3").Value = "Volkswagen" + [char]10 + "Golf"
$ws.Range("A4").Value = "Toyota" + [char]10 + "Corolla"
$ws.Range("A5").Value = "Skoda" + [char]10 + "Octavia"
$ws.Range("A6").Value = "BMW" + [char]10 + "3" + [char]10 + "Series"
$ws.Range("A7").Value = "Hyundai" + [char]10 + "Solaris"

# --- Updated priority-vector values ---------------------------------
# Keep these as text (leading "'" forces text, matching the original
# cells which are stored as text rather than numbers).
$ws.Range("H4").Value = "'0.228"
$ws.Range("H6").Value = "'0.362"
